# Sample Template.docx improvement:
#   Split the paragraph that holds "<para1>" + line breaks + "Regards" into three
#   paragraphs, inserting two new placeholder lines ("<style> " and "<table>") that
#   mirror the existing "<Signature>" placeholder pattern used elsewhere in the doc.
#
# Before:
#   <para1>[br][br][br]Regards           (all one paragraph)
#
# After:
#   <para1>                              (paragraph 1)
#   <style> [br]                         (paragraph 2 - new)
#   <table>[br][br]Regards               (paragraph 3)

$d = $word.ActiveDocument

# Locate the "<para1>" placeholder text; collapse the range to just after it.
$find = $d.Content
$found = $find.Find.Execute("<para1>", $false, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '<para1>' placeholder in the document."
}
$afterPara1 = $find.End

# Split the paragraph right after "<para1>" -- this separates it from the trailing
# line-break runs, becoming its own paragraph.
$splitPoint = $d.Range($afterPara1, $afterPara1)
$splitPoint.InsertParagraphAfter()

# The very next paragraph now starts with the original run of manual line breaks
# followed by "Regards". Drop the first of those line breaks (it is being replaced
# by the new paragraph mark above).
$restPara = $d.Paragraphs.Item($d.Paragraphs.Item(1).Index)  # placeholder, reset below
$restPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq ($afterPara1 + 1)) {
        $restPara = $cand
        break
    }
}
if ($null -eq $restPara) {
    # Fallback: the paragraph immediately following the one that now ends at
    # (afterPara1 + 1) -- locate by scanning paragraph starts.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Start -ge $afterPara1) {
            $restPara = $cand
            break
        }
    }
}

$firstBreak = $d.Range($restPara.Range.Start, $restPara.Range.Start + 1)
$firstBreak.Text = ""

# Anchor sits right before the (now second, originally third) manual line break /
# "Regards" run -- this run keeps its original identity, so new content inserted
# immediately before it is kept in its own separate run/paragraph instead of being
# merged into it.
$restPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq ($afterPara1 + 1)) {
        $restPara = $cand
        break
    }
}
$anchor = $restPara.Range.Start

# Insert the new "<table>" placeholder text right before the remaining content.
$p = $d.Range($anchor, $anchor)
$p.InsertBefore("<table>")

# Insert a new manual line break before "<table>" -- this becomes the end of the
# new "<style> " paragraph.
$p = $d.Range($anchor, $anchor)
$p.InsertBefore("`v")

# Insert the "<style> " placeholder text as three separate runs, matching the
# existing "<Signature>" placeholder's run pattern ("<" / name / ">").
$p = $d.Range($anchor, $anchor)
$p.InsertBefore("> ")
$p = $d.Range($anchor, $anchor)
$p.InsertBefore("style")
$p = $d.Range($anchor, $anchor)
$p.InsertBefore("<")

# Finally, split the new content into its own paragraph: "<style> " + the new
# break become paragraph 2, "<table>" + the remaining original content become
# paragraph 3. The boundary sits right after "<style> " (8 characters) and the new
# break (1 character).
$splitOffset = $anchor + "<style> ".Length + 1
$splitPoint2 = $d.Range($splitOffset, $splitOffset)
$splitPoint2.InsertParagraphBefore()
